$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells for the season record columns ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the formatting of the existing last header cell (AB1) onto the
# three new header cells so they match the rest of row 1 (bold, centered,
# bordered header style).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the season record (Wins/Losses/Ties) for every player row ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 96   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 66   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
